$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5 (shifts the existing numbered rows down by one)
# and add the new Italian phrase "e mezzanotte" ("and midnight") as a new shared string.
$ws.Rows.Item(5).Insert()
$ws.Range("A5").Value = "e mezzanotte"

# Normalize formatting: switch column A's cells from the old "loud" (bold/europe
# themed) text style to a plain Arial 10 text format, matching the rest of the
# workbook's original font instead of the louder default.
$plain = $wb.Styles.Add("PlainText")
$plain.Font.Name = "Arial"
$plain.Font.Size = 10
$ws.Columns.Item(1).Style = "PlainText"
$ws.Range("A1:A64").NumberFormat = "@"

# Move the active selection to C58, as recorded for the sheet view.
$ws.Range("C58").Select()
